# Insert a new data row at row 52 (shifting existing rows 52:158 down to 53:159)
# and populate it with the new weekly record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("52:52").Insert()

$ws.Range("A52").Value = 7
$ws.Range("B52").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C52").Value = "Ñuble"
$ws.Range("D52").Value = 44533
$ws.Range("E52").Value = 16
$ws.Range("F52").Value = 100112003
$ws.Range("G52").Value = "Ajo"
$ws.Range("H52").Value = "Chino"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 60
$ws.Range("K52").Value = 20000
$ws.Range("L52").Value = 21000
$ws.Range("M52").Value = 20500
$ws.Range("N52").Value = "$/caja 10 kilos"
$ws.Range("O52").Value = "China"
$ws.Range("P52").Value = 2050
$ws.Range("Q52").Value = 10
$ws.Range("R52").Value = "Hortaliza"
